# edit.ps1 - Applies the cryptos.xlsx price/volume/ranking refresh
# described in the commit "Updated cryptos list ... with GitHub Actions".
#
# For cells whose new text looks like a plain number (e.g. "9.56"), Excel's
# COM layer auto-converts a plain Range.Value assignment into a numeric cell
# (and can even drop trailing zeros, e.g. "94.40" -> 94.4). The source workbook
# stores these as text (t="inlineStr"), so for those cells we briefly force the
# cell to Text format, assign the value, then ClearFormats() to drop back to the
# default (unstyled) cell format while keeping the stored value as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "70.268.03"
$ws.Range("E2").Value = "  +0.35%  "
# Row 3
$ws.Range("D3").Value = "3.559.52"
$ws.Range("E3").Value = "  +0.39%  "
# Row 4
$ws.Range("E4").Value = "  -0.08%  "
# Row 5
Set-TextValue "D5" "607.66"
$ws.Range("E5").Value = "  +3.39%  "
# Row 6
Set-TextValue "D6" "186.25"
$ws.Range("E6").Value = "  +0.66%  "
# Row 7
$ws.Range("D7").Value = "3.553.87"
$ws.Range("E7").Value = "  +0.52%  "
# Row 8
Set-TextValue "D8" "0.618"
$ws.Range("E8").Value = "  +0.40%  "
# Row 9
$ws.Range("E9").Value = "  -0.01%  "
# Row 10
Set-TextValue "D10" "0.216"
$ws.Range("E10").Value = "  +9.18%  "
# Row 11
Set-TextValue "D11" "0.645"
$ws.Range("E11").Value = "  +0.11%  "
# Row 12
$ws.Range("E12").Value = "  -0.65%  "
# Row 13
$ws.Range("E13").Value = "  +1.43%  "
# Row 14
Set-TextValue "D14" "9.56"
$ws.Range("E14").Value = "  +0.73%  "
# Row 15
$ws.Range("D15").Value = "4.124.21"
$ws.Range("E15").Value = "  +0.42%  "
# Row 16
$ws.Range("D16").Value = "70.328.57"
$ws.Range("E16").Value = "  +0.50%  "
# Row 17
$ws.Range("D17").Value = "3.567.58"
$ws.Range("E17").Value = "  +1.48%  "
# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D18" "12.72"
$ws.Range("E18").Value = "  +2.25%  "
# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "19.03"
$ws.Range("E19").Value = "  -1.76%  "
# Row 20
Set-TextValue "D20" "579.67"
$ws.Range("E20").Value = "  +7.35%  "
# Row 21
$ws.Range("E21").Value = "  +0.48%  "
# Row 22
$ws.Range("E22").Value = "  -1.94%  "
# Row 23
Set-TextValue "D23" "17.38"
$ws.Range("E23").Value = "  -3.41%  "
# Row 24
$ws.Range("E24").Value = "  +2.66%  "
# Row 25
Set-TextValue "D25" "4.89"
$ws.Range("E25").Value = "  +0.58%  "
# Row 26
Set-TextValue "D26" "94.40"
$ws.Range("E26").Value = "  -1.31%  "
# Row 27
Set-TextValue "D27" "2.94"
$ws.Range("E27").Value = "  -1.85%  "
# Row 28
Set-TextValue "D28" "10.96"
$ws.Range("E28").Value = "  -2.17%  "
# Row 29
Set-TextValue "D29" "9.40"
$ws.Range("E29").Value = "  +2.96%  "
# Row 30
$ws.Range("E30").Value = "  +0.47%  "
# Row 31
Set-TextValue "D31" "7.06"
$ws.Range("E31").Value = "  -3.18%  "
# Row 32
$ws.Range("E32").Value = "  -1.57%  "
# Row 33
Set-TextValue "D33" "0.114"
$ws.Range("E33").Value = "  +1.42%  "
# Row 34
Set-TextValue "D34" "63.57"
$ws.Range("E34").Value = "  -1.49%  "
# Row 35
Set-TextValue "D35" "3.69"
$ws.Range("E35").Value = "  +19.18%  "
# Row 36
$ws.Range("E36").Value = "  -1.20%  "
# Row 37
Set-TextValue "D37" "527.60"
$ws.Range("E37").Value = "  -3.65%  "
# Row 38
$ws.Range("E38").Value = "  -1.90%  "
# Row 39
$ws.Range("E39").Value = "  -0.10%  "
# Row 40
Set-TextValue "D40" "37.47"
$ws.Range("E40").Value = "  -2.13%  "
# Row 41
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0786"
$ws.Range("E41").Value = "  +2.74%  "
# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.546.82"
$ws.Range("E42").Value = "  +5.84%  "
# Row 43
Set-TextValue "D43" "3.56"
$ws.Range("E43").Value = "  +4.55%  "
# Row 44
$ws.Range("E44").Value = "  +1.52%  "
# Row 45
$ws.Range("E45").Value = "  +4.20%  "
# Row 46
$ws.Range("E46").Value = "  -1.34%  "
# Row 47
Set-TextValue "D47" "3.42"
$ws.Range("E47").Value = "  -4.48%  "
# Row 48
Set-TextValue "D48" "0.141"
$ws.Range("E48").Value = "  +3.51%  "
# Row 49
Set-TextValue "D49" "9.19"
$ws.Range("E49").Value = "  +0.37%  "
# Row 50
$ws.Range("E50").Value = "  +0.30%  "
# Row 51
Set-TextValue "D51" "135.64"
$ws.Range("E51").Value = "  -1.17%  "
